$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 116.09524
$ws.Range("J9").Value = 155.42857
$ws.Range("L9").Value = 155.42857
$ws.Range("N9").Value = -493.42857
$ws.Range("H62").Value = 145841390
$ws.Range("I62").Value = 62511850
$ws.Range("K62").Value = 62511850
$ws.Range("M62").Value = -62511226
$ws.Range("H65").Value = 145841390
$ws.Range("I65").Value = 62511850
$ws.Range("K65").Value = 312559250
$ws.Range("M65").Value = -312556130
$ws.Range("I98").Value = 6897871
$ws.Range("K98").Value = 6897871
$ws.Range("M98").Value = -6896373
$ws.Range("H106").Value = 58826210
$ws.Range("I106").Value = 142860140
$ws.Range("J106").Value = 2450
$ws.Range("K106").Value = 142860140
$ws.Range("L106").Value = 2450
$ws.Range("M106").Value = -142859509
$ws.Range("N106").Value = -3712
$ws.Range("I122").Value = 6897871
$ws.Range("K122").Value = 20693613
$ws.Range("M122").Value = -20691163
$ws.Range("H137").Value = 40130040
$ws.Range("I137").Value = 69445310
$ws.Range("K137").Value = 208335930
$ws.Range("M137").Value = -208333380

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7581063.5
$ws.Range("I32").Value = 5311.741
$ws.Range("K32").Value = 5311.741
$ws.Range("M32").Value = -5024.741
$ws.Range("H61").Value = 4171249.8
$ws.Range("I61").Value = 3631252.5
$ws.Range("J61").Value = 5349425.5
$ws.Range("K61").Value = 3631252.5
$ws.Range("L61").Value = 5349425.5
$ws.Range("M61").Value = -3631040.5
$ws.Range("N61").Value = -5349849.5
$ws.Range("H122").Value = 2425.3845
$ws.Range("I122").Value = 2360.3157
$ws.Range("J122").Value = 2602
$ws.Range("K122").Value = 7080.9471
$ws.Range("L122").Value = 7806
$ws.Range("M122").Value = -4630.9471
$ws.Range("N122").Value = -12706
$ws.Range("H132").Value = 19896910
$ws.Range("I132").Value = 20068068
$ws.Range("J132").Value = 19659924
$ws.Range("K132").Value = 60204204
$ws.Range("L132").Value = 58979772
$ws.Range("M132").Value = -60201674
$ws.Range("N132").Value = -58984832
$ws.Range("H136").Value = 4171249.8
$ws.Range("I136").Value = 3631252.5
$ws.Range("J136").Value = 5349425.5
$ws.Range("K136").Value = 10893757.5
$ws.Range("L136").Value = 16048276.5
$ws.Range("M136").Value = -10891207.5
$ws.Range("N136").Value = -16053376.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22059798
$ws.Range("I134").Value = 33334210
$ws.Range("J134").Value = 3269110.8
$ws.Range("K134").Value = 100002630
$ws.Range("L134").Value = 9807332.399999999
$ws.Range("M134").Value = -100000095
$ws.Range("N134").Value = -9812402.399999999
$ws.Range("H140").Value = 29890
$ws.Range("J140").Value = 29890
$ws.Range("L140").Value = 29890
$ws.Range("N140").Value = -40250

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13515045
$ws.Range("I31").Value = 25001102
$ws.Range("J31").Value = 2036.1471
$ws.Range("K31").Value = 25001102
$ws.Range("L31").Value = 2036.1471
$ws.Range("M31").Value = -25000807
$ws.Range("N31").Value = -2626.1471
$ws.Range("H34").Value = 13515045
$ws.Range("I34").Value = 25001102
$ws.Range("J34").Value = 2036.1471
$ws.Range("K34").Value = 25001102
$ws.Range("L34").Value = 2036.1471
$ws.Range("M34").Value = -25000900
$ws.Range("N34").Value = -2440.1471
$ws.Range("H58").Value = 981898.0600000001
$ws.Range("I58").Value = 1401370.8
$ws.Range("J58").Value = 3128.4443
$ws.Range("K58").Value = 1401370.8
$ws.Range("L58").Value = 3128.4443
$ws.Range("M58").Value = -1401167.8
$ws.Range("N58").Value = -3534.4443
$ws.Range("H122").Value = 11012.444
$ws.Range("I122").Value = 13516
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 40548
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -38098
$ws.Range("N122").Value = -11650
$ws.Range("H132").Value = 2085484.4
$ws.Range("I132").Value = 3847215.2
$ws.Range("K132").Value = 11541645.6
$ws.Range("M132").Value = -11539115.6
$ws.Range("H134").Value = 1178501
$ws.Range("I134").Value = 1522.0714
$ws.Range("K134").Value = 4566.2142
$ws.Range("M134").Value = -2031.2142
$ws.Range("H136").Value = 981898.0600000001
$ws.Range("I136").Value = 1401370.8
$ws.Range("J136").Value = 3128.4443
$ws.Range("K136").Value = 4204112.4
$ws.Range("L136").Value = 9385.332900000001
$ws.Range("M136").Value = -4201562.4
$ws.Range("N136").Value = -14485.3329

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2751.018
$ws.Range("I68").Value = 577.625
$ws.Range("J68").Value = 4433.645
$ws.Range("K68").Value = 1732.875
$ws.Range("L68").Value = 13300.935
$ws.Range("M68").Value = -921.875
$ws.Range("N68").Value = -14922.935
$ws.Range("H71").Value = 2751.018
$ws.Range("I71").Value = 577.625
$ws.Range("J71").Value = 4433.645
$ws.Range("K71").Value = 5198.625
$ws.Range("L71").Value = 39902.80500000001
$ws.Range("M71").Value = -1142.625
$ws.Range("N71").Value = -48014.80500000001
$ws.Range("H107").Value = 394930.4
$ws.Range("I107").Value = 625599.6
$ws.Range("J107").Value = 870.5
$ws.Range("K107").Value = 1876798.8
$ws.Range("L107").Value = 2611.5
$ws.Range("M107").Value = -1874878.8
$ws.Range("N107").Value = -6451.5
$ws.Range("H131").Value = 8351027.5
$ws.Range("I131").Value = 500000000
$ws.Range("J131").Value = 17994.236
$ws.Range("K131").Value = 1500000000
$ws.Range("L131").Value = 53982.708
$ws.Range("M131").Value = -1499994960
$ws.Range("N131").Value = -64062.708
$ws.Range("H132").Value = 3263.2222
$ws.Range("I132").Value = 2321.3333
$ws.Range("J132").Value = 3734.1667
$ws.Range("K132").Value = 20891.9997
$ws.Range("L132").Value = 33607.5003
$ws.Range("M132").Value = -18361.9997
$ws.Range("N132").Value = -38667.5003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 15788
$ws.Range("J49").Value = 15788
$ws.Range("L49").Value = 15788
$ws.Range("N49").Value = -16156
$ws.Range("H122").Value = 5954397.5
$ws.Range("I122").Value = 1858.1818
$ws.Range("J122").Value = 27780376
$ws.Range("K122").Value = 5574.5454
$ws.Range("L122").Value = 83341128
$ws.Range("M122").Value = -3124.5454
$ws.Range("N122").Value = -83346028
$ws.Range("H132").Value = 19324914
$ws.Range("I132").Value = 16429252
$ws.Range("J132").Value = 23378840
$ws.Range("K132").Value = 49287756
$ws.Range("L132").Value = 70136520
$ws.Range("M132").Value = -49285226
$ws.Range("N132").Value = -70141580

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4055688.2
$ws.Range("I132").Value = 5136158.5
$ws.Range("J132").Value = 3924.5
$ws.Range("K132").Value = 15408475.5
$ws.Range("L132").Value = 11773.5
$ws.Range("M132").Value = -15405945.5
$ws.Range("N132").Value = -16833.5
$ws.Range("H136").Value = 28826566
$ws.Range("I136").Value = 9821472
$ws.Range("K136").Value = 29464416
$ws.Range("M136").Value = -29461866

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 59980
$ws.Range("J131").Value = 59980
$ws.Range("L131").Value = 59980
$ws.Range("N131").Value = -70060
$ws.Range("H132").Value = 290290.78
$ws.Range("I132").Value = 397931.75
$ws.Range("J132").Value = 3248.1667
$ws.Range("K132").Value = 1193795.25
$ws.Range("L132").Value = 9744.500100000001
$ws.Range("M132").Value = -1191265.25
$ws.Range("N132").Value = -14804.5001
$ws.Range("H136").Value = 5069.636
$ws.Range("I136").Value = 3815.76
$ws.Range("J136").Value = 6719.4736
$ws.Range("K136").Value = 11447.28
$ws.Range("L136").Value = 20158.4208
$ws.Range("M136").Value = -8897.280000000001
$ws.Range("N136").Value = -25258.4208
